$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top, pushing all existing rows (old 1-18)
# down to (new 2-19).
$ws.Rows("1:1").Insert()

# Put the new intro sentence into the freshly inserted A1.
$ws.Range("A1").Value = "This is a rough demo of the client / server protocol."

# Restore a sane view: scroll back to the top-left and select G8, as in the
# saved workbook (previously the view was scrolled to A9 with D21 selected).
$ws.Range("G8").Select()
